# Re-run of the signup/login integration test captured another round of
# tokens + ids for the three users (daniel5f, Jorge2525, mario35). Update the
# "token" (C) and "id" (D) columns on rows 2-4 with the freshly generated
# values from the latest test run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - daniel5f
$ws.Range("C2").Value = "eyJhbGciOiJIUzI1NiIsInR5cCI6IkpXVCJ9.eyJ1c2VyTmFtZSI6ImRhbmllbDVmIiwicGFzc3dvcmQiOiJBejI1Mjg4QCIsImlhdCI6MTcwMTMwODk1OX0.NZwYWfyRDPUY3cBFpURzeqqiLx4JL2-wxV6Ogq8qVEQ"
$ws.Range("D2").Value = "488e3f8c-8bc8-4047-9e46-56ef8a7d3339"

# Row 3 - Jorge2525
$ws.Range("C3").Value = "eyJhbGciOiJIUzI1NiIsInR5cCI6IkpXVCJ9.eyJ1c2VyTmFtZSI6IkpvcmdlMjUyNSIsInBhc3N3b3JkIjoiYXNUMzU2NDQ0QCIsImlhdCI6MTcwMTMwODk2MH0.GNR_HM9RLhboVsRtFFh8zMByru_0JH4U_qkkWUsWSkk"
$ws.Range("D3").Value = "3f70a5cd-7551-4bb7-aa81-93ad2fcbecac"

# Row 4 - mario35
$ws.Range("C4").Value = "eyJhbGciOiJIUzI1NiIsInR5cCI6IkpXVCJ9.eyJ1c2VyTmFtZSI6Im1hcmlvMzUiLCJwYXNzd29yZCI6Im1BcmlvdXVnQDMiLCJpYXQiOjE3MDEzMDg5NjF9.8e6br2l_Jnmu8PfBAfLg3meBKRAnOgGkwRTkwvV1CxU"
$ws.Range("D4").Value = "aae53cfd-232f-48c5-928f-eb2463ccfaaa"
